# Fixed medication export issue AB#12255
#
# The "Brand" column (C) in the medication report template was bound to
# the same placeholder as the "Generic" column, so the export showed the
# generic name twice instead of the brand name. Point C2/C3 at dedicated
# brand placeholders for each record row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = "{d.records[i].brand}"
$ws.Range("C3").Value2 = "{d.records[i+1].brand}"

# Leave the selection on the last cell touched, matching the author's
# final editing position.
$ws.Range("C3").Select()
